# Auto-generated edit script applying numeric updates to the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# across several leve-profit worksheets, per the scheduled-runner
# refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4922.25
$ws.Range("I18").Value = 4922.25
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 4922.25
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -4638.25
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 2175.5417
$ws.Range("I32").Value = 2064
$ws.Range("J32").Value = 2204.8948
$ws.Range("K32").Value = 2064
$ws.Range("L32").Value = 2204.8948
$ws.Range("M32").Value = -1738
$ws.Range("N32").Value = -2856.8948
$ws.Range("H86").Value = 2895.3157
$ws.Range("I86").Value = 3213.8333
$ws.Range("J86").Value = 2349.2856
$ws.Range("K86").Value = 3213.8333
$ws.Range("L86").Value = 2349.2856
$ws.Range("M86").Value = -2090.8333
$ws.Range("N86").Value = -4595.2856
$ws.Range("H89").Value = 2895.3157
$ws.Range("I89").Value = 3213.8333
$ws.Range("J89").Value = 2349.2856
$ws.Range("K89").Value = 16069.1665
$ws.Range("L89").Value = 11746.428
$ws.Range("M89").Value = -10453.1665
$ws.Range("N89").Value = -22978.428
$ws.Range("H98").Value = 8160.88
$ws.Range("I98").Value = 334.90475
$ws.Range("J98").Value = 49247.25
$ws.Range("K98").Value = 334.90475
$ws.Range("L98").Value = 49247.25
$ws.Range("M98").Value = 1163.09525
$ws.Range("N98").Value = -52243.25
$ws.Range("H103").Value = 696.5
$ws.Range("I103").Value = 650
$ws.Range("J103").Value = 761.6
$ws.Range("K103").Value = 1950
$ws.Range("L103").Value = 2284.8
$ws.Range("M103").Value = -1364
$ws.Range("N103").Value = -3456.8
$ws.Range("H112").Value = 2114.9412
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 2246.7144
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 6740.1432
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -8956.143199999999
$ws.Range("H122").Value = 8160.88
$ws.Range("I122").Value = 334.90475
$ws.Range("J122").Value = 49247.25
$ws.Range("K122").Value = 1004.71425
$ws.Range("L122").Value = 147741.75
$ws.Range("M122").Value = 1445.28575
$ws.Range("N122").Value = -152641.75
$ws.Range("H137").Value = 1668.1975
$ws.Range("I137").Value = 1057.491
$ws.Range("J137").Value = 2960.077
$ws.Range("K137").Value = 3172.473
$ws.Range("L137").Value = 8880.231
$ws.Range("M137").Value = -622.473
$ws.Range("N137").Value = -13980.231
$ws.Range("H138").Value = 3296.9453
$ws.Range("I138").Value = 2847.8965
$ws.Range("J138").Value = 3592.9092
$ws.Range("K138").Value = 8543.6895
$ws.Range("L138").Value = 10778.7276
$ws.Range("M138").Value = -3403.6895
$ws.Range("N138").Value = -21058.7276
$ws.Range("H141").Value = 3534.4211
$ws.Range("I141").Value = 2872.8057
$ws.Range("J141").Value = 15443.5
$ws.Range("K141").Value = 8618.417099999999
$ws.Range("L141").Value = 46330.5
$ws.Range("M141").Value = -3438.417099999999
$ws.Range("N141").Value = -56690.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2661.88
$ws.Range("I32").Value = 1972.7285
$ws.Range("J32").Value = 12310
$ws.Range("K32").Value = 1972.7285
$ws.Range("L32").Value = 12310
$ws.Range("M32").Value = -1685.7285
$ws.Range("N32").Value = -12884
$ws.Range("H74").Value = 2311.7715
$ws.Range("I74").Value = 1731.3334
$ws.Range("J74").Value = 2512.6924
$ws.Range("K74").Value = 1731.3334
$ws.Range("L74").Value = 2512.6924
$ws.Range("M74").Value = -857.3334
$ws.Range("N74").Value = -4260.6924
$ws.Range("H77").Value = 2311.7715
$ws.Range("I77").Value = 1731.3334
$ws.Range("J77").Value = 2512.6924
$ws.Range("K77").Value = 8656.666999999999
$ws.Range("L77").Value = 12563.462
$ws.Range("M77").Value = -4288.666999999999
$ws.Range("N77").Value = -21299.462
$ws.Range("H110").Value = 1282.2727
$ws.Range("I110").Value = 1210.6
$ws.Range("J110").Value = 1999
$ws.Range("K110").Value = 1210.6
$ws.Range("L110").Value = 1999
$ws.Range("M110").Value = 834.4000000000001
$ws.Range("N110").Value = -6089
$ws.Range("H122").Value = 2360.9119
$ws.Range("I122").Value = 1402.9584
$ws.Range("J122").Value = 4660
$ws.Range("K122").Value = 4208.8752
$ws.Range("L122").Value = 13980
$ws.Range("M122").Value = -1758.8752
$ws.Range("N122").Value = -18880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 237.88889
$ws.Range("I80").Value = 99.5
$ws.Range("J80").Value = 277.42856
$ws.Range("K80").Value = 99.5
$ws.Range("L80").Value = 277.42856
$ws.Range("M80").Value = 898.5
$ws.Range("N80").Value = -2273.42856
$ws.Range("H83").Value = 237.88889
$ws.Range("I83").Value = 99.5
$ws.Range("J83").Value = 277.42856
$ws.Range("K83").Value = 497.5
$ws.Range("L83").Value = 1387.1428
$ws.Range("M83").Value = 4494.5
$ws.Range("N83").Value = -11371.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3464.4375
$ws.Range("I31").Value = 2050.8823
$ws.Range("J31").Value = 5066.467
$ws.Range("K31").Value = 2050.8823
$ws.Range("L31").Value = 5066.467
$ws.Range("M31").Value = -1755.8823
$ws.Range("N31").Value = -5656.467
$ws.Range("H34").Value = 3464.4375
$ws.Range("I34").Value = 2050.8823
$ws.Range("J34").Value = 5066.467
$ws.Range("K34").Value = 2050.8823
$ws.Range("L34").Value = 5066.467
$ws.Range("M34").Value = -1848.8823
$ws.Range("N34").Value = -5470.467
$ws.Range("H58").Value = 2267.0527
$ws.Range("I58").Value = 1501
$ws.Range("J58").Value = 2956.5
$ws.Range("K58").Value = 1501
$ws.Range("L58").Value = 2956.5
$ws.Range("M58").Value = -1298
$ws.Range("N58").Value = -3362.5
$ws.Range("H105").Value = 808.6667
$ws.Range("I105").Value = 813.375
$ws.Range("J105").Value = 771
$ws.Range("K105").Value = 813.375
$ws.Range("L105").Value = 771
$ws.Range("M105").Value = 933.625
$ws.Range("N105").Value = -4265
$ws.Range("H132").Value = 3274.2246
$ws.Range("I132").Value = 2969.3416
$ws.Range("J132").Value = 4836.75
$ws.Range("K132").Value = 8908.024800000001
$ws.Range("L132").Value = 14510.25
$ws.Range("M132").Value = -6378.024800000001
$ws.Range("N132").Value = -19570.25
$ws.Range("H134").Value = 2908.5454
$ws.Range("I134").Value = 2251.1538
$ws.Range("J134").Value = 8036.2
$ws.Range("K134").Value = 6753.4614
$ws.Range("L134").Value = 24108.6
$ws.Range("M134").Value = -4218.4614
$ws.Range("N134").Value = -29178.6
$ws.Range("H136").Value = 2267.0527
$ws.Range("I136").Value = 1501
$ws.Range("J136").Value = 2956.5
$ws.Range("K136").Value = 4503
$ws.Range("L136").Value = 8869.5
$ws.Range("M136").Value = -1953
$ws.Range("N136").Value = -13969.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 822.6
$ws.Range("I5").Value = 851.5
$ws.Range("J5").Value = 803.3333
$ws.Range("K5").Value = 2554.5
$ws.Range("L5").Value = 2409.9999
$ws.Range("M5").Value = -2442.5
$ws.Range("N5").Value = -2633.9999
$ws.Range("H56").Value = 6774.7676
$ws.Range("I56").Value = 6774.7676
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 6774.7676
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -6244.7676
$ws.Range("H135").Value = 822.6
$ws.Range("I135").Value = 851.5
$ws.Range("J135").Value = 803.3333
$ws.Range("K135").Value = 7663.5
$ws.Range("L135").Value = 7229.9997
$ws.Range("M135").Value = -5128.5
$ws.Range("N135").Value = -12299.9997
$ws.Range("H136").Value = 3451.8333
$ws.Range("I136").Value = 3451.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10355.4999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5255.499899999999
$ws.Range("H138").Value = 50013704
$ws.Range("I138").Value = 90913096
$ws.Range("J138").Value = 25555
$ws.Range("K138").Value = 272739288
$ws.Range("L138").Value = 76665
$ws.Range("M138").Value = -272734148
$ws.Range("N138").Value = -86945
$ws.Range("H141").Value = 24248.12
$ws.Range("I141").Value = 6305.7617
$ws.Range("J141").Value = 42190.477
$ws.Range("K141").Value = 18917.2851
$ws.Range("L141").Value = 126571.431
$ws.Range("M141").Value = -13737.2851
$ws.Range("N141").Value = -136931.431

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7653.1665
$ws.Range("I70").Value = 7333.6665
$ws.Range("J70").Value = 8611.666999999999
$ws.Range("K70").Value = 7333.6665
$ws.Range("L70").Value = 8611.666999999999
$ws.Range("M70").Value = -7063.6665
$ws.Range("N70").Value = -9151.666999999999
$ws.Range("H73").Value = 7653.1665
$ws.Range("I73").Value = 7333.6665
$ws.Range("J73").Value = 8611.666999999999
$ws.Range("K73").Value = 7333.6665
$ws.Range("L73").Value = 8611.666999999999
$ws.Range("M73").Value = -6397.6665
$ws.Range("N73").Value = -10483.667
$ws.Range("H97").Value = 1203.5714
$ws.Range("I97").Value = 506.64285
$ws.Range("J97").Value = 2597.4285
$ws.Range("K97").Value = 506.64285
$ws.Range("L97").Value = 2597.4285
$ws.Range("M97").Value = -10.64285000000001
$ws.Range("N97").Value = -3589.4285
$ws.Range("H126").Value = 6449
$ws.Range("I126").Value = 9596.666999999999
$ws.Range("J126").Value = 5399.778
$ws.Range("K126").Value = 28790.001
$ws.Range("L126").Value = 16199.334
$ws.Range("M126").Value = -26320.001
$ws.Range("N126").Value = -21139.334
$ws.Range("H132").Value = 1737.7693
$ws.Range("I132").Value = 1634.3055
$ws.Range("J132").Value = 2979.3333
$ws.Range("K132").Value = 4902.916499999999
$ws.Range("L132").Value = 8937.999899999999
$ws.Range("M132").Value = -2372.916499999999
$ws.Range("N132").Value = -13997.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1179.7273
$ws.Range("I22").Value = 899.75
$ws.Range("J22").Value = 1339.7142
$ws.Range("K22").Value = 899.75
$ws.Range("L22").Value = 1339.7142
$ws.Range("M22").Value = -604.75
$ws.Range("N22").Value = -1929.7142
$ws.Range("H27").Value = 1179.7273
$ws.Range("I27").Value = 899.75
$ws.Range("J27").Value = 1339.7142
$ws.Range("K27").Value = 899.75
$ws.Range("L27").Value = 1339.7142
$ws.Range("M27").Value = -792.75
$ws.Range("N27").Value = -1553.7142
$ws.Range("H132").Value = 1990.0526
$ws.Range("I132").Value = 1658.8918
$ws.Range("J132").Value = 2602.7
$ws.Range("K132").Value = 4976.6754
$ws.Range("L132").Value = 7808.099999999999
$ws.Range("M132").Value = -2446.6754
$ws.Range("N132").Value = -12868.1
$ws.Range("H136").Value = 1595.0754
$ws.Range("I136").Value = 1179.925
$ws.Range("J136").Value = 2872.4614
$ws.Range("K136").Value = 3539.775
$ws.Range("L136").Value = 8617.3842
$ws.Range("M136").Value = -989.7749999999996
$ws.Range("N136").Value = -13717.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2030.2142
$ws.Range("I126").Value = 1569
$ws.Range("J126").Value = 3183.25
$ws.Range("K126").Value = 4707
$ws.Range("L126").Value = 9549.75
$ws.Range("M126").Value = -2237
$ws.Range("N126").Value = -14489.75
$ws.Range("H132").Value = 2563.7222
$ws.Range("I132").Value = 2608.0334
$ws.Range("J132").Value = 2342.1667
$ws.Range("K132").Value = 7824.100199999999
$ws.Range("L132").Value = 7026.500100000001
$ws.Range("M132").Value = -5294.100199999999
$ws.Range("N132").Value = -12086.5001
$ws.Range("H136").Value = 1795.8049
$ws.Range("I136").Value = 1866.3871
$ws.Range("J136").Value = 1577
$ws.Range("K136").Value = 5599.1613
$ws.Range("L136").Value = 4731
$ws.Range("M136").Value = -3049.1613
$ws.Range("N136").Value = -9831
